$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the previously empty/placeholder predicted price: convert the
# hard-coded "TimeTaken in Hours" value in C2 into a real formula that
# derives it from the minutes value in B2.
$ws.Range("C2").Formula = "=B2/60"
